$wb = $excel.ActiveWorkbook

# --- Sheet "VerifyCSVForExistingVersion" (was empty H column, now Pass/Fail) ---
$ws3 = $wb.Worksheets.Item("VerifyCSVForExistingVersion")

$passRows = @(5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,23)
foreach ($r in $passRows) {
    $ws3.Range("H$r").Value = "Pass"
}

$failRows = @(20,21,22)
foreach ($r in $failRows) {
    $ws3.Range("H$r").Value = "Fail"
}

# --- Sheet "Cases_RealTimeSpine" (previously all Pass, now cleared) ---
$ws7 = $wb.Worksheets.Item("Cases_RealTimeSpine")
for ($r = 5; $r -le 80; $r++) {
    $ws7.Range("H$r").ClearContents()
}

# --- View / selection state ---
# Sheet7 loses its "tabSelected" flag and gets a new selection
[void]$ws7.Activate()
[void]$ws7.Range("F83").Select()

# Sheet3 becomes the active / selected sheet with a new selection
[void]$ws3.Activate()
[void]$ws3.Range("D20").Select()
